# Add two new worksheets ("reserved" and "probab") to the workbook, both
# built as copies of the existing "states" sheet (same layout/styles),
# positioned immediately before "states". Then tweak one cell on each new
# sheet, and leave "states" as the selected/active sheet - matching the
# target revision of inst/test_cases/test_states.xlsx.

$wb = $excel.ActiveWorkbook

$statesSheet = $wb.Worksheets.Item("states")

# Copy "states" twice; each Copy($statesSheet) inserts the new copy
# immediately before $statesSheet, so doing this twice yields, in order:
#   ... dup_names, <copy #1>, <copy #2>, states
# Excel names the copies "states (2)" / "states (2) (2)".
$statesSheet.Copy($statesSheet)
$statesSheet.Copy($statesSheet)

$wb.Worksheets.Item("states (2) (2)").Name = "reserved"
$wb.Worksheets.Item("states (2)").Name = "probab"

# "reserved": the name/id column now reads "cycle" instead of "relapse",
# typed fresh (no inherited cell style).
$reserved = $wb.Worksheets.Item("reserved")
$reserved.Activate() | Out-Null
$reserved.Range("A4").Style = "Normal"
$reserved.Range("A4").Value = "cycle"
$reserved.Range("C13").Select() | Out-Null

# "probab": initial_probability for the first row becomes 0.5.
$probab = $wb.Worksheets.Item("probab")
$probab.Activate() | Out-Null
$probab.Range("D2").Value = 0.5
$probab.Range("D2").Select() | Out-Null

# The "wrong_name" sheet picks up explicit best-fit widths for columns
# F:H (it already holds data there, just without a stored <col> width).
$wrongName = $wb.Worksheets.Item("wrong_name")
$wrongName.Columns("F").ColumnWidth = 14.830729166666666
$wrongName.Columns("G").ColumnWidth = 14.498697916666666
$wrongName.Columns("H").ColumnWidth = 18.830729166666668

# "states" stays/becomes the visible, active sheet. Re-fetch by name:
# $statesSheet was captured before the two inserts-before-it above, and a
# worksheet variable tracks its slot, not the renamed/shifted original.
$wb.Worksheets.Item("states").Activate() | Out-Null
